# Updates cryptos list: prices and 1h volume percentages for multiple rows
# (also swaps the Kaspa/VeChain row order)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = '66.204.32'
$ws.Range("E2").Value = '  -1.39%  '
$ws.Range("D3").Value = '3.209.43'
$ws.Range("E3").Value = '  -0.26%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''606.38'
$ws.Range("E5").Value = '  +0.29%  '
$ws.Range("D6").Value = '''154.86'
$ws.Range("E6").Value = '  -1.92%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '3.211.05'
$ws.Range("E8").Value = '  -0.24%  '
$ws.Range("D9").Value = '''0.545'
$ws.Range("E9").Value = '  -2.54%  '
$ws.Range("D11").Value = '''5.70'
$ws.Range("E11").Value = '  -3.87%  '
$ws.Range("D12").Value = '''0.499'
$ws.Range("E12").Value = '  -4.16%  '
$ws.Range("D13").Value = '''0.0000266'
$ws.Range("E13").Value = '  -1.83%  '
$ws.Range("D14").Value = '''38.14'
$ws.Range("E14").Value = '  -3.62%  '
$ws.Range("D15").Value = '3.733.94'
$ws.Range("E15").Value = '  -0.28%  '
$ws.Range("D16").Value = '66.284.46'
$ws.Range("E16").Value = '  -1.04%  '
$ws.Range("D17").Value = '3.208.90'
$ws.Range("E17").Value = '  -0.33%  '
$ws.Range("D18").Value = '''7.23'
$ws.Range("E18").Value = '  -4.13%  '
$ws.Range("E19").Value = '  +0.74%  '
$ws.Range("D20").Value = '''505.24'
$ws.Range("E20").Value = '  -4.18%  '
$ws.Range("D21").Value = '''15.19'
$ws.Range("E21").Value = '  -2.63%  '
$ws.Range("E22").Value = '  -3.08%  '
$ws.Range("E23").Value = '  -3.89%  '
$ws.Range("D24").Value = '''14.48'
$ws.Range("E24").Value = '  -4.39%  '
$ws.Range("D25").Value = '''84.78'
$ws.Range("E25").Value = '  -1.31%  '
$ws.Range("D26").Value = '''0.152'
$ws.Range("E26").Value = '  +66.78%  '
$ws.Range("D27").Value = '''1.00'
$ws.Range("E27").Value = '  -0.03%  '
$ws.Range("E28").Value = '  -1.61%  '
$ws.Range("D29").Value = '''8.94'
$ws.Range("E29").Value = '  -4.62%  '
$ws.Range("E30").Value = '  -3.40%  '
$ws.Range("E31").Value = '  -2.24%  '
$ws.Range("E32").Value = '  -4.04%  '
$ws.Range("D33").Value = '''28.16'
$ws.Range("E33").Value = '  -1.23%  '
$ws.Range("E34").Value = '  +0.09%  '
$ws.Range("D35").Value = '''1.16'
$ws.Range("E35").Value = '  -5.72%  '
$ws.Range("E36").Value = '  -4.25%  '
$ws.Range("D37").Value = '''55.32'
$ws.Range("E37").Value = '  +0.20%  '
$ws.Range("D38").Value = '''497.85'
$ws.Range("E38").Value = '  -5.98%  '
$ws.Range("D39").Value = '0.0₃0759'
$ws.Range("E39").Value = '  +10.05%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = '''0.0416'
$ws.Range("E40").Value = '  -3.23%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").Value = '''0.129'
$ws.Range("E41").Value = '  +0.79%  '
$ws.Range("D42").Value = '''2.99'
$ws.Range("E42").Value = '  +1.89%  '
$ws.Range("D43").Value = '''8.69'
$ws.Range("E44").Value = '  -3.82%  '
$ws.Range("D45").Value = '2.912.78'
$ws.Range("E45").Value = '  -0.10%  '
$ws.Range("D46").Value = '''2.43'
$ws.Range("E46").Value = '  -2.23%  '
$ws.Range("D47").Value = '''27.88'
$ws.Range("E47").Value = '  -3.32%  '
$ws.Range("E48").Value = '  +0.10%  '
$ws.Range("E50").Value = '  -1.34%  '
$ws.Range("D51").Value = '''121.59'
$ws.Range("E51").Value = '  -0.63%  '

# Clear the quote-prefix formatting introduced above so cell styling
# matches the original (unstyled) text cells.
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D51").Style = "Normal"
